# "fix(publipostage): Try to solve Excel emoji problem"
#
# The "statut" column (A) used four book emoji as status markers:
#   (red book)    -> now the text  "-3"
#   (blue book)   -> now the emoji "warning sign"
#   (orange book) -> now the text  "+3"
#   (green book)  -> now the emoji "check mark"
#
# Every cell that used one of the old emoji gets its text replaced with the
# corresponding new marker below - nothing else on the sheet changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldRed    = "📕"
$oldBlue   = "📘"
$oldOrange = "📙"
$oldGreen  = "📗"

$newRed    = "-3"
$newBlue   = "⚠️"
$newOrange = "+3"
$newGreen  = "✅"

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow  = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()

    if ($old -eq $oldRed -or $old -eq $oldBlue -or $old -eq $oldOrange -or $old -eq $oldGreen) {
        if ($old -eq $oldRed) {
            $newText = $newRed
        } elseif ($old -eq $oldBlue) {
            $newText = $newBlue
        } elseif ($old -eq $oldOrange) {
            $newText = $newOrange
        } else {
            $newText = $newGreen
        }

        # "-3"/"+3" read as numeric literals to the value-setter, so force a
        # text number format while writing them, then drop the cell back to
        # the sheet's normal (unstyled) look - these cells carry no
        # explicit style either before or after the edit.
        $cell.NumberFormat = "@"
        $cell.Value = $newText
        $cell.Style = "Normal"
    }
}
